# Update "Forecast Comparison" sheet with corrected forecast output:
#  - insert a new "Week_Start_Date" column after "Week" (shifts ASIN.. right by one column)
#  - normalise the Week labels from "W01".."W16" to "W1".."W16"
#  - populate the new Week_Start_Date column with the Monday-of-week dates (stored as text)
#  - store is_holiday_week as a boolean instead of a number

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new blank column at B; everything from the old B (ASIN) onward
# shifts one column to the right (B->C, C->D, ... I->J).
$ws.Columns("B:B").Insert()

# Make sure the new column is treated as plain text so the ISO date strings
# are not reinterpreted as date serial numbers.
$ws.Range("B1:B17").NumberFormat = "@"

$ws.Range("B1").Value = "Week_Start_Date"

$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $weekStartDates[$i]
    $ws.Range("A$row").Value = "W" + ($i + 1)
    $ws.Range("J$row").Value = $false
}
